# Replace embedded line breaks with spaces in several vaccine brand/manufacturer
# name cells across the two Influenza worksheets. This collapses a formerly
# duplicated shared string ("Fluvirin Preservative-free") automatically.

$wb = $excel.ActiveWorkbook

$pedFlu = $wb.Worksheets.Item("Pediatric Influenza Vaccine ")
$pedFlu.Range("B3").Value = "Fluzone Pediatric dose No Preservative"
$pedFlu.Range("B6").Value = "Fluarix Preservative-Free"
$pedFlu.Range("B9").Value = "FluMist No Preservative"
$pedFlu.Range("B10").Value = "Afluria No Preservative"
$pedFlu.Range("H10").Value = "Merck (CSL product)"

$adultFlu = $wb.Worksheets.Item("Adult Influenza Vaccine ")
$adultFlu.Range("B5").Value = "Agriflu No Preservative"
$adultFlu.Range("B7").Value = "Fluvirin Preservative-free"
$adultFlu.Range("B8").Value = "Fluarix Preservative-free"
$adultFlu.Range("B10").Value = "Flumist No Preservative"
